$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row for new columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the style of an existing header cell (AC1) to the new headers
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats

# Fill data rows 2-51 with the season record values
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 30).Value = 62   # AD
    $ws.Cells.Item($r, 31).Value = 99   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
